# Insert 4 new rows before row 840 (shifting old rows 840-876 down to 844-880)
# and populate the new rows with the new weekly price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A840:A843").EntireRow.Insert()

# Copy the number format used by the date column (D) from an existing data row
# so the newly inserted date cells render the same way.
$dateFormat = $ws.Range("D2").NumberFormat

$newRows = @(
    @{ Row = 840; D = 45075; L = "Especial"; M = 100; N = 35000; O = 35000; P = 35000; Q = "$/bandeja 10 kilos"; R = "Perú"; S = 3500; T = 10 },
    @{ Row = 841; D = 45075; L = "Especial"; M = 100; N = 4500;  O = 4500;  P = 4500;  Q = "$/kilo (en caja de 17 kilos)"; R = "Provincia de Quillota"; S = 4500; T = 1 },
    @{ Row = 842; D = 45075; L = "Primera";  M = 180; N = 28000; O = 30000; P = 29111; Q = "$/bandeja 10 kilos"; R = "Perú"; S = 2911; T = 10 },
    @{ Row = 843; D = 45075; L = "Primera";  M = 140; N = 3000;  O = 4000;  P = 3571;  Q = "$/kilo (en caja de 17 kilos)"; R = "Provincia de Quillota"; S = 3571; T = 1 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 7
    $ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($row, 3).Value = "Ñuble"

    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 4).NumberFormat = $dateFormat

    $ws.Cells.Item($row, 5).Value = 16
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100106
    $ws.Cells.Item($row, 8).Value = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value = 100106002
    $ws.Cells.Item($row, 10).Value = "Palta"
    $ws.Cells.Item($row, 11).Value = "Hass"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
